$wb = $excel.ActiveWorkbook

# --- Work on "Management Interface (Mobiel)" sheet (4th tab) ---
$ws = $wb.Worksheets.Item(4)

# Rows 14-17: rename the browser-only entries to "PC / <browser>"
$ws.Range("B14").Value = "PC / Chrome"
$ws.Range("B15").Value = "PC / Safari"
$ws.Range("B16").Value = "PC / Firefox"
$ws.Range("B17").Value = "PC / Internet Explorer"

# Rows 18-21: new test points for mobile browsers
$ws.Range("B18").Value = "Android / AndroidBrowser (Internet)"
$ws.Range("C18").Value = "Closed"
$ws.Range("E18").Value = 41624
$ws.Range("F18").Value = "Hoog"

$ws.Range("B19").Value = "Android / Chrome"
$ws.Range("C19").Value = "Closed"
$ws.Range("E19").Value = 41624
$ws.Range("F19").Value = "Hoog"

$ws.Range("B20").Value = "iOS / Safari"
$ws.Range("C20").Value = "Closed"
$ws.Range("E20").Value = 41624
$ws.Range("F20").Value = "Hoog"

$ws.Range("B21").Value = "iOS / Chrome"
$ws.Range("C21").Value = "Closed"
$ws.Range("E21").Value = 41624
$ws.Range("F21").Value = "Hoog"

# Give the new rows the same look as the existing, similarly-wrapped rows
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 30

# Make this sheet the active tab/selection, matching the saved view state
$ws.Activate()
$ws.Range("B3").Select()
